# TokenIteratorFieldRewriterSplit: rewrite the " m:'doc.html'.fromHTMLURI() "
# Word field (fldChar begin/instrText.../fldChar end) into plain-text runs
# that spell out the M2Doc token delimiters "{" ... "}", keeping the
# _GoBack bookmark in place between "doc.html" and "'.fromHTMLURI()".

$d = $word.ActiveDocument

# Locate the paragraph that contains the field (rather than hard-coding an
# index) so the script is resilient to minor structural differences.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Fields.Count -gt 0) {
        $targetPara = $candidate
        break
    }
}

if ($targetPara -eq $null) {
    throw "Could not find a paragraph containing a field"
}

# Build the replacement paragraph's WordprocessingML: the field's begin
# fldChar and the leading/trailing " " instrText runs disappear; every
# remaining instrText run becomes a plain <w:t> run, and the bookmark stays
# exactly where it was, between the "doc.html" run and the "'.fromHTMLURI()"
# run.
$openTick = [char]39

$paraXml = '<w:p w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F">' +
    '<w:r><w:t>{</w:t></w:r>' +
    '<w:r><w:t>m</w:t></w:r>' +
    '<w:r><w:t>:</w:t></w:r>' +
    '<w:r><w:t>' + $openTick + '</w:t></w:r>' +
    '<w:r><w:t>doc.html</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t>' + $openTick + '.fromHTMLURI()</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">}</w:t></w:r>' +
    '</w:p>'

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' + $paraXml + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$targetPara.Range.InsertXML($packageXml)
